# "Removal of unnecessary config after 4G section"
#
# - Region (B2) changes from APAC to EMEA
# - Backup Link's "4G+Cellular" flag (B24) is switched from False to True
# - Active cell selection on the sheet moves to G14

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "EMEA"
$ws.Range("B24").Value = $true

$ws.Range("G14").Select() | Out-Null
